# Update countries & provincias Spain
# Applies the 13-Jul-2020 12:06 -> 13:23 data refresh to the "Pais" sheet:
#  - Refreshes several per-country metrics (Casos totales, Nuevos casos,
#    Casos activos, Recuperados, Muertes hoy, Muertes).
#  - A handful of countries overtook their neighbours in the "Casos
#    totales" ranking, so their rows swap places (the row keeps its
#    position in the sheet, but the two/three country names trade
#    places and the newly-updated country's figures move up with its
#    name).
#  - Updates the "Datos actualizados a ..." timestamp banner in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Timestamp banner
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 13 de Julio de 2020 a las 13:23"

# ---------------------------------------------------------------------
# 2) Country-name swaps caused by re-ranking (column A)
# ---------------------------------------------------------------------
$ws.Range("A50").Value  = "Suiza"
$ws.Range("A51").Value  = "Barein"

$ws.Range("A73").Value  = "Kenia"
$ws.Range("A74").Value  = "Sudan"

$ws.Range("A95").Value  = "Madagascar"
$ws.Range("A96").Value  = "Republica de Yibuti"
$ws.Range("A97").Value  = "Luxemburgo"

$ws.Range("A124").Value = "Cabo Verde"
$ws.Range("A125").Value = "Sierra Leona"

$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"

# ---------------------------------------------------------------------
# 3) Refreshed numeric figures (columns B:H)
#    B = Casos totales, C = Nuevos casos, D = Casos activos,
#    E = Recuperados, F = Casos criticos, G = Muertes hoy, H = Muertes
# ---------------------------------------------------------------------

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 3414105
$ws.Range("C4").Value = 110
$ws.Range("D4").Value = 1517560
$ws.Range("E4").Value = 1758758
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 137787

# Row 6 - India
$ws.Range("B6").Value = 879902
$ws.Range("C6").Value = 436
$ws.Range("D6").Value = 554907
$ws.Range("E6").Value = 301795

# Row 14 - Iran
$ws.Range("B14").Value = 259652
$ws.Range("C14").Value = 2349
$ws.Range("D14").Value = 222539
$ws.Range("E14").Value = 24081
$ws.Range("G14").Value = 203
$ws.Range("H14").Value = 13032

# Row 19 - Alemania
$ws.Range("B19").Value = 199998
$ws.Range("C19").Value = 48
$ws.Range("E19").Value = 6263

# Row 50 - Suiza (was Barein's row; now shows Suiza's fresh figures)
$ws.Range("B50").Value = 32946
$ws.Range("C50").Value = 63
$ws.Range("D50").Value = 29600
$ws.Range("E50").Value = 1378
$ws.Range("H50").Value = 1968

# Row 51 - Barein (was Suiza's row; now shows Barein's carried-over figures)
$ws.Range("B51").Value = 32941
$ws.Range("D51").Value = 28425
$ws.Range("E51").Value = 4408
$ws.Range("H51").Value = 108

# Row 64 - Nepal
$ws.Range("B64").Value = 16945
$ws.Range("C64").Value = 144
$ws.Range("D64").Value = 10294
$ws.Range("E64").Value = 6613

# Row 73 - Kenia (was Sudan's row; now shows Kenia's fresh figures)
$ws.Range("B73").Value = 10294
$ws.Range("C73").Value = 189
$ws.Range("D73").Value = 2946
$ws.Range("E73").Value = 7151
$ws.Range("G73").Value = 12
$ws.Range("H73").Value = 197

# Row 74 - Sudan (was Kenia's row; now shows Sudan's carried-over figures)
$ws.Range("B74").Value = 10250
$ws.Range("D74").Value = 5341
$ws.Range("E74").Value = 4259
$ws.Range("H74").Value = 650

# Row 80 - Senegal
$ws.Range("B80").Value = 8198
$ws.Range("C80").Value = 63
$ws.Range("D80").Value = 5514
$ws.Range("E80").Value = 2534
$ws.Range("G80").Value = 2
$ws.Range("H80").Value = 150

# Row 95 - Madagascar (fresh figures, moved up two spots)
$ws.Range("B95").Value = 5080
$ws.Range("C95").Value = 213
$ws.Range("D95").Value = 2494
$ws.Range("E95").Value = 2549
$ws.Range("G95").Value = 2
$ws.Range("H95").Value = 37

# Row 96 - Republica de Yibuti (carried-over figures, shifted down one)
$ws.Range("B96").Value = 4972
$ws.Range("D96").Value = 4712
$ws.Range("E96").Value = 204
$ws.Range("H96").Value = 56

# Row 97 - Luxemburgo (carried-over figures, shifted down one)
$ws.Range("B97").Value = 4925
$ws.Range("D97").Value = 4086
$ws.Range("E97").Value = 728
$ws.Range("H97").Value = 111

# Row 124 - Cabo Verde (was Sierra Leona's row; now shows Cabo Verde's fresh figures)
$ws.Range("B124").Value = 1698
$ws.Range("C124").Value = 75
$ws.Range("D124").Value = 748
$ws.Range("E124").Value = 931
$ws.Range("H124").Value = 19

# Row 125 - Sierra Leona (was Cabo Verde's row; now shows Sierra Leona's carried-over figures)
$ws.Range("B125").Value = 1635
$ws.Range("D125").Value = 1154
$ws.Range("E125").Value = 418
$ws.Range("H125").Value = 63
